$wb = $excel.ActiveWorkbook

# --- Sheet "Metadata" (sheet1) ---
$ws1 = $wb.Worksheets.Item("Metadata")

# Remove the duplicate "Contact" row (row 11); row 10's "Contact" row will be
# repurposed into the new "Jurisdiction" row below.
$ws1.Rows.Item(11).Delete()

# Version 5.0.0 -> 6.0.0
$ws1.Cells.Item(3, 2).Value = "6.0.0"

# Date updated
$ws1.Cells.Item(8, 2).Value = "2022-01-21T20:46:54+00:00"

# Publisher now has a value
$ws1.Cells.Item(9, 2).Value = "Alvearie Team"

# Former "Contact" row (row 10) becomes "Jurisdiction" / "United States of America"
$ws1.Cells.Item(10, 1).Value = "Jurisdiction"
$ws1.Cells.Item(10, 2).Value = "United States of America"

# Case Sensitive now has value "true" (force text so it isn't stored as a
# Boolean TRUE value)
$ws1.Cells.Item(14, 2).Value = "'true"

# --- Sheet "Concepts" (sheet2) ---
$ws2 = $wb.Worksheets.Item("Concepts")

# Move the "adv-imputed-pcp" row (currently the last row, row 6) so that it
# becomes the first concept row (row 2), right after the header.
$ws2.Rows.Item(2).Insert()
$ws2.Range("A7:D7").Copy($ws2.Range("A2:D2"))
$ws2.Rows.Item(7).Delete()
